$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27, shifting existing rows 27:44 down to 28:45.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new weekly data point.
# (Same static descriptive fields as the surrounding rows; new price data.)
$ws.Cells.Item(27, 1).Value2 = 1
$ws.Cells.Item(27, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value2 = 44438
$ws.Cells.Item(27, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 5).Value2 = 15
$ws.Cells.Item(27, 6).Value2 = 100112040
$ws.Cells.Item(27, 7).Value2 = "Cilantro"
$ws.Cells.Item(27, 8).Value2 = "Sin especificar"
$ws.Cells.Item(27, 9).Value2 = "Primera"
$ws.Cells.Item(27, 10).Value2 = 300
$ws.Cells.Item(27, 11).Value2 = 1000
$ws.Cells.Item(27, 12).Value2 = 1200
$ws.Cells.Item(27, 13).Value2 = 1100
$ws.Cells.Item(27, 14).Value2 = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(27, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value2 = 550
$ws.Cells.Item(27, 17).Value2 = 2
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"
